# Applies the crypto price/volume refresh described in the commit:
# "Updated cryptos list on Wed Oct 23 08:46:53 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.773.50"
$ws.Range("E2").Value = "  -0.52%  "

# Row 3
$ws.Range("D3").Value = "2.591.19"
$ws.Range("E3").Value = "  -1.78%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.81"
$ws.Range("E5").Value = "  -1.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.03"
$ws.Range("E6").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  -3.05%  "

# Row 9
$ws.Range("D9").Value = "2.600.41"
$ws.Range("E9").Value = "  -1.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  -3.11%  "

# Row 11
$ws.Range("E11").Value = "  +0.38%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("E12").Value = "  +0.21%  "

# Row 13
$ws.Range("E13").Value = "  -1.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.04"
$ws.Range("E14").Value = "  -2.38%  "

# Row 15
$ws.Range("D15").Value = "3.074.37"
$ws.Range("E15").Value = "  -1.38%  "

# Row 16
$ws.Range("E16").Value = "  -2.18%  "

# Row 17
$ws.Range("D17").Value = "66.828.76"
$ws.Range("E17").Value = "  -0.55%  "

# Row 18
$ws.Range("D18").Value = "2.602.43"
$ws.Range("E18").Value = "  -1.24%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.56"
$ws.Range("E19").Value = "  -4.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.76"
$ws.Range("E20").Value = "  -4.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.37"
$ws.Range("E21").Value = "  -2.42%  "

# Row 22
$ws.Range("E22").Value = "  -2.61%  "

# Row 23
$ws.Range("E23").Value = "  -3.45%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.03%  "

# Row 25
$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.33"
$ws.Range("E25").Value = "  -6.05%  "

# Row 26
$ws.Range("E26").Value = "  -4.76%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.04"
$ws.Range("E27").Value = "  -2.41%  "

# Row 28
$ws.Range("D28").Value = "2.735.70"
$ws.Range("E28").Value = "  -1.26%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.21%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0983"
$ws.Range("E30").Value = "  -3.31%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "538.41"
$ws.Range("E31").Value = "  -2.81%  "

# Row 32
$ws.Range("E32").Value = "  +1.58%  "

# Row 33
$ws.Range("E33").Value = "  -3.82%  "

# Row 34
$ws.Range("E34").Value = "  -3.30%  "

# Row 35
$ws.Range("E35").Value = "  -0.80%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.12%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.47"
$ws.Range("E37").Value = "  -4.23%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.53"
$ws.Range("E38").Value = "  -0.14%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.79"
$ws.Range("E39").Value = "  -2.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.361"
$ws.Range("E40").Value = "  -2.44%  "

# Row 41
$ws.Range("E41").Value = "  +2.04%  "

# Row 42
$ws.Range("E42").Value = "  -0.95%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.09"
$ws.Range("E43").Value = "  -3.30%  "

# Row 44
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("E45").Value = "  -4.52%  "

# Row 46
$ws.Range("D46").Value = "0.0₆0291"
$ws.Range("E46").Value = "  -1.89%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.87"
$ws.Range("E47").Value = "  -2.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.570"
$ws.Range("E48").Value = "  -3.70%  "

# Row 49
$ws.Range("E49").Value = "  -3.00%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.69"
$ws.Range("E50").Value = "  -1.84%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0762"
$ws.Range("E51").Value = "  -1.95%  "

